# Apply the edits described in the commit "Edits to progress report".
#
# Paragraph (Section 1 summary of the Colin Bailey article):
#   1) Insert "Keogh River " before "steelhead, cutthroat trout, and"
#   2) Change "and also body sizes" -> "and body sizes" (drops the stray
#      "also" along with the now-stale grammar-check proofErr markers
#      that wrapped it)
#   3) Insert "on the Keogh " before "(see section 2)."

$d = $word.ActiveDocument

# 1) "...on how steelhead, cutthroat trout, and..." -> "...on how Keogh River steelhead, cutthroat trout, and..."
$d.Content.Find.Execute(
    "steelhead, cutthroat trout, and", $true, $false, $false, $false, $false,
    $true, 1, $false, "Keogh River steelhead, cutthroat trout, and", 2) | Out-Null

# 2) "...species, densities, and also body sizes..." -> "...species, densities, and body sizes..."
#    Use a range that spans real text on both sides of the old "and also" /
#    proofErr markers so the stale grammar-check markers around "and also"
#    are cleared away rather than left dangling around the replacement.
$rng = $d.Content
$rng.Find.Execute("densities, and also body sizes") | Out-Null
$rng.Text = "densities, and body sizes"

# 3) "...broader population-level consequences (see section 2)." ->
#    "...broader population-level consequences on the Keogh (see section 2)."
$d.Content.Find.Execute(
    "consequences (see section 2).", $true, $false, $false, $false, $false,
    $true, 1, $false, "consequences on the Keogh (see section 2).", 2) | Out-Null
